# Edit script: Update Demand_Tech_Unit_Cost with improved cost assumptions
# and adjust workbook view/selection state.

$wb = $excel.ActiveWorkbook

# --- Sheet activation / view state -----------------------------------
# Target state: "Input_Params" becomes the selected tab (was
# "Demand_Tech_Characteristics"), and "Demand_Tech_Unit_Cost" loses its
# stale B1 scroll position / E7 selection (reset back towards A1).
$wsCost = $wb.Worksheets.Item("Demand_Tech_Unit_Cost")
$wsCost.Activate()
$wsCost.Range("A1").Select()

$wsInput = $wb.Worksheets.Item("Input_Params")
$wsInput.Activate()

$ws = $wsCost

# --- Demand_Tech_Unit_Cost data updates ---
$ws.Range("D2").Value = 16500
$ws.Range("E2").Formula = "=D2-(D2*0.02)"
$ws.Range("F2").Formula = "=E2-(E2*0.02)"
$ws.Range("G2").Formula = "=F2-(F2*0.02)"
$ws.Range("H2").Formula = "=G2-(G2*0.02)"
$ws.Range("I2").Formula = "=H2-(H2*0.02)"

$ws.Range("D3").Value = 5900
$ws.Range("E3").Formula = "=D3-(D3*0.015)"
$ws.Range("F3").Formula = "=E3-(E3*0.015)"
$ws.Range("G3").Formula = "=F3-(F3*0.015)"
$ws.Range("H3").Formula = "=G3-(G3*0.015)"
$ws.Range("I3").Formula = "=H3-(H3*0.015)"

$ws.Range("D4").Value = 4
$ws.Range("E4").Formula = "=D4-(D4*0.03)"
$ws.Range("F4").Formula = "=E4-(E4*0.03)"
$ws.Range("G4").Formula = "=F4-(F4*0.03)"
$ws.Range("H4").Formula = "=G4-(G4*0.03)"
$ws.Range("I4").Formula = "=H4-(H4*0.03)"

$ws.Range("D5").Value = 130
$ws.Range("E5").Value = 128
$ws.Range("F5").Value = 126
$ws.Range("G5").Value = 124
$ws.Range("H5").Value = 122
$ws.Range("I5").Value = 120

$ws.Range("D6").Value = 1495
$ws.Range("E6").Formula = "=D6-(D6*0.0375)"
$ws.Range("F6").Formula = "=E6-(E6*0.0375)"
$ws.Range("G6").Formula = "=F6-(F6*0.0375)"
$ws.Range("H6").Formula = "=G6-(G6*0.0375)"
$ws.Range("I6").Formula = "=H6-(H6*0.0375)"

$ws.Range("D7").Value = 1700
$ws.Range("E7").Value = 1670
$ws.Range("F7").Value = 1640
$ws.Range("G7").Value = 1610
$ws.Range("H7").Value = 1580
$ws.Range("I7").Value = 1550

$ws.Range("D8").Value = 1360
$ws.Range("E8").Formula = "=D8-(D8*0.095)"
$ws.Range("F8").Formula = "=E8-(E8*0.095)"
$ws.Range("G8").Formula = "=F8-(F8*0.095)"
$ws.Range("H8").Formula = "=G8-(G8*0.095)"
$ws.Range("I8").Formula = "=H8-(H8*0.095)"

$ws.Range("D9").Value = 2500
$ws.Range("E9").Value = 2450
$ws.Range("F9").Value = 2400
$ws.Range("G9").Value = 2350
$ws.Range("H9").Value = 2300
$ws.Range("I9").Value = 2250

$ws.Range("D10").Value = 3500
$ws.Range("E10").Value = 3430
$ws.Range("F10").Value = 3360
$ws.Range("G10").Value = 3290
$ws.Range("H10").Value = 3220
$ws.Range("I10").Value = 3150

$ws.Range("D11").Value = 6000
$ws.Range("E11").Formula = "=D11-(D11*0.02)"
$ws.Range("F11").Formula = "=E11-(E11*0.02)"
$ws.Range("G11").Formula = "=F11-(F11*0.02)"
$ws.Range("H11").Formula = "=G11-(G11*0.02)"
$ws.Range("I11").Formula = "=H11-(H11*0.02)"

$ws.Range("D12").Value = 6500
$ws.Range("E12").Value = 6435
$ws.Range("F12").Value = 6370
$ws.Range("G12").Value = 6305
$ws.Range("H12").Value = 6240
$ws.Range("I12").Value = 6175

$ws.Range("D13").Value = 24000
$ws.Range("E13").Value = 23500
$ws.Range("F13").Value = 23000
$ws.Range("G13").Value = 22500
$ws.Range("H13").Value = 22000
$ws.Range("I13").Value = 21500

$ws.Range("D14").Value = 22000
$ws.Range("E14").Value = 21600
$ws.Range("F14").Value = 21200
$ws.Range("G14").Value = 20800
$ws.Range("H14").Value = 20400
$ws.Range("I14").Value = 20000

$ws.Range("D15").Value = 35000
$ws.Range("E15").Value = 34500
$ws.Range("F15").Value = 34000
$ws.Range("G15").Value = 33500
$ws.Range("H15").Value = 33000
$ws.Range("I15").Value = 32500

$ws.Range("D16").Value = 2800
$ws.Range("E16").Value = 2750
$ws.Range("F16").Value = 2700
$ws.Range("G16").Value = 2650
$ws.Range("H16").Value = 2600
$ws.Range("I16").Value = 2550

$ws.Range("D17").Value = 1200
$ws.Range("E17").Value = 1180
$ws.Range("F17").Value = 1160
$ws.Range("G17").Value = 1140
$ws.Range("H17").Value = 1120
$ws.Range("I17").Value = 1100

$ws.Range("D18").Value = 140
$ws.Range("E18").Value = 138
$ws.Range("F18").Value = 136
$ws.Range("G18").Value = 134
$ws.Range("H18").Value = 132
$ws.Range("I18").Value = 130

$ws.Range("D19").Value = 450
$ws.Range("E19").Value = 445
$ws.Range("F19").Value = 440
$ws.Range("G19").Value = 435
$ws.Range("H19").Value = 430
$ws.Range("I19").Value = 425

$ws.Range("D20").Value = 2200
$ws.Range("E20").Value = 2150
$ws.Range("F20").Value = 2100
$ws.Range("G20").Value = 2050
$ws.Range("H20").Value = 2000
$ws.Range("I20").Value = 1950

$ws.Range("D21").Value = 6500
$ws.Range("E21").Value = 6400
$ws.Range("F21").Value = 6300
$ws.Range("G21").Value = 6200
$ws.Range("H21").Value = 6100
$ws.Range("I21").Value = 6000

$ws.Range("D22").Value = 75
$ws.Range("E22").Value = 74
$ws.Range("F22").Value = 73
$ws.Range("G22").Value = 72
$ws.Range("H22").Value = 71
$ws.Range("I22").Value = 70


# --- Number formatting --------------------------------------------------
# All of the updated data cells (D2:I22) pick up a new "0.00" number
# format (numFmtId 2) that didn't exist in the original style table.
$ws.Range("D2:I22").NumberFormat = "0.00"

Write-Host "Demand_Tech_Unit_Cost updated"
